$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.840.05'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '2.291.05'
$ws.Range('E3').Value = '  -0.16%  '
$__style = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = $__style
$ws.Range('E4').Value = '  -0.07%  '
$__style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.46'
$ws.Range('D5').Style = $__style
$ws.Range('E5').Value = '  +0.21%  '
$__style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.90'
$ws.Range('D6').Style = $__style
$ws.Range('E6').Value = '  +2.40%  '
$ws.Range('E7').Value = '  -0.33%  '
$ws.Range('E8').Value = '  +0.02%  '
$__style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.509'
$ws.Range('D9').Style = $__style
$ws.Range('E9').Value = '  +3.09%  '
$__style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.96'
$ws.Range('D10').Style = $__style
$ws.Range('E10').Value = '  +7.93%  '
$__style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0789'
$ws.Range('D11').Style = $__style
$ws.Range('E11').Value = '  -0.64%  '
$ws.Range('E12').Value = '  +1.29%  '
$__style = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.26'
$ws.Range('D13').Style = $__style
$ws.Range('E13').Value = '  +8.67%  '
$__style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.91'
$ws.Range('D14').Style = $__style
$ws.Range('E14').Value = '  +2.29%  '
$ws.Range('D15').Value = '2.649.45'
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('D16').Value = '2.295.88'
$ws.Range('E16').Value = '  +0.94%  '
$__style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.798'
$ws.Range('D17').Style = $__style
$ws.Range('E17').Value = '  -0.94%  '
$ws.Range('D18').Value = '42.750.16'
$ws.Range('E18').Value = '  +0.06%  '
$__style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.41'
$ws.Range('D19').Style = $__style
$ws.Range('E19').Value = '  +7.99%  '
$ws.Range('E20').Value = '  +0.39%  '
$__style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.10'
$ws.Range('D21').Style = $__style
$ws.Range('E21').Value = '  +1.56%  '
$ws.Range('E22').Value = '  +0.69%  '
$__style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.44'
$ws.Range('D23').Style = $__style
$ws.Range('E23').Value = '  +0.01%  '
$__style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.23'
$ws.Range('D24').Style = $__style
$ws.Range('E24').Value = '  +11.95%  '
$__style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = $__style
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('E26').Value = '  -0.28%  '
$__style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.90'
$ws.Range('D27').Style = $__style
$ws.Range('E27').Value = '  +2.77%  '
$ws.Range('E28').Value = '  +14.89%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$__style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '167.20'
$ws.Range('D29').Style = $__style
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$__style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.34'
$ws.Range('D30').Style = $__style
$ws.Range('E30').Value = '  +1.90%  '
$__style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.12'
$ws.Range('D31').Style = $__style
$ws.Range('E31').Value = '  +0.55%  '
$__style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.999'
$ws.Range('D32').Style = $__style
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('E33').Value = '  +1.62%  '
$__style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.63'
$ws.Range('D34').Style = $__style
$ws.Range('E34').Value = '  +5.43%  '
$ws.Range('E35').Value = '  -0.98%  '
$ws.Range('E36').Value = '  +0.97%  '
$__style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0686'
$ws.Range('D37').Style = $__style
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$__style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.78'
$ws.Range('D39').Style = $__style
$ws.Range('E39').Value = '  +2.33%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$__style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.101'
$ws.Range('D40').Style = $__style
$ws.Range('E40').Value = '  +0.20%  '
$ws.Range('E41').Value = '  +0.57%  '
$__style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.31'
$ws.Range('D42').Style = $__style
$ws.Range('E42').Value = '  -0.97%  '
$ws.Range('E43').Value = '  +4.25%  '
$ws.Range('D44').Value = '1.968.84'
$ws.Range('E44').Value = '  -1.04%  '
$__style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.14'
$ws.Range('D45').Style = $__style
$ws.Range('E45').Value = '  +3.57%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$__style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.88'
$ws.Range('D46').Style = $__style
$ws.Range('E46').Value = '  +1.66%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$__style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.50'
$ws.Range('D47').Style = $__style
$ws.Range('E47').Value = '  +0.47%  '
$__style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.14'
$ws.Range('D48').Style = $__style
$ws.Range('E48').Value = '  +4.72%  '
$ws.Range('E49').Value = '  +3.69%  '
$ws.Range('D50').Value = '2.517.13'
$ws.Range('E50').Value = '  -0.13%  '
$__style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '70.60'
$ws.Range('D51').Style = $__style
$ws.Range('E51').Value = '  +1.38%  '
